# feat: add 2022-Q3 data
#
# 1) Insert a new row into the "总计" (summary) sheet for the 2022-Q3 quarter.
# 2) Insert a new worksheet "2022-Q3" (placed right after "总计") holding the
#    per-fund detail rows for that quarter.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) "总计" summary sheet: insert a new row right under the header with the
#    2022-Q3 totals; everything below shifts down automatically.
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item(1)

$summary.Rows.Item(2).Insert()

# column A is a simple running index (0,1,2,...) styled like the cell below it
$summary.Cells.Item(3,1).Copy()
$summary.Cells.Item(2,1).PasteSpecial(-4122)
$summary.Range("B2:D2").Style = "Normal"

$summary.Cells.Item(2,1).Value = 0
$summary.Cells.Item(2,2).Value = "2022-Q3"
$summary.Cells.Item(2,3).Value = 8
$summary.Cells.Item(2,4).Value = 3.25

# ---------------------------------------------------------------------------
# 2) New "2022-Q3" worksheet: duplicate the existing "2022-Q2" sheet (so it
#    inherits the same header/column styling) and drop it right after "总计".
# ---------------------------------------------------------------------------
$template = $wb.Worksheets.Item(2)
$template.Copy($null, $summary)

$q3 = $wb.Worksheets.Item(2)
$q3.Name = "2022-Q3"

# the template only carried one data row; make room for all eight
$q3.Rows.Item(3).Resize(7).Insert()

# re-apply the index-column style (row-insert can blend in stray formatting)
$q3.Range("A2").Copy()
$q3.Range("A3:A9").PasteSpecial(-4122)

$funds = @(
    @(0, "001532", "华安文体健康主题灵活配置混合A", "42.79", "87.37", "2.63", "1.1254", 5),
    @(1, "002350", "华安安华灵活配置混合A",         "32.74", "87.55", "2.35", "0.7694", 6),
    @(2, "014207", "华安产业精选混合A",             "26.65", "82.10", "1.94", "0.5170", 7),
    @(3, "014208", "华安产业精选混合C",             "22.90", "82.10", "1.94", "0.4443", 7),
    @(4, "009970", "财通内需增长12个月定期开放混合", "9.38",  "56.38", "2.28", "0.2139", 8),
    @(5, "013116", "华安文体健康主题灵活配置混合C", "4.16",  "87.37", "2.63", "0.1094", 5),
    @(6, "080005", "长盛量化红利混合",               "1.89",  "61.68", "3.28", "0.0620", 2),
    @(7, "016183", "华安安华灵活配置混合C",         "0.45",  "87.55", "2.35", "0.0106", 6)
)

for ($i = 0; $i -lt $funds.Count; $i++) {
    $r = $i + 2
    $fund = $funds[$i]
    $q3.Cells.Item($r, 1).Value = $fund[0]
    $q3.Cells.Item($r, 2).Value = "'" + $fund[1]
    $q3.Cells.Item($r, 3).Value = $fund[2]
    $q3.Cells.Item($r, 4).Value = "'" + $fund[3]
    $q3.Cells.Item($r, 5).Value = "'" + $fund[4]
    $q3.Cells.Item($r, 6).Value = "'" + $fund[5]
    $q3.Cells.Item($r, 7).Value = "'" + $fund[6]
    $q3.Cells.Item($r, 8).Value = $fund[7]
}

# clear the stray text-format styling the quote-prefixed writes above leave
# behind, matching the plain (unstyled) data cells used elsewhere
$q3.Range("B2:H9").Style = "Normal"
